$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.459612070389937
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 3.781711156805759

$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.1575252929769615
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 3.781711156805759
